$wb = $excel.ActiveWorkbook

# ALC row 12
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 93.14286
$ws.Range("I12").Value = 92
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 92
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = 78

# ALC row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# ALC row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

# ALC row 64
$ws.Range("H64").Value = 2970
$ws.Range("I64").Value = 2744.4443
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 2744.4443
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -2496.4443
$ws.Range("N64").Value = -5496

# ALC row 67
$ws.Range("H67").Value = 2970
$ws.Range("I67").Value = 2744.4443
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 2744.4443
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -1886.4443
$ws.Range("N67").Value = -6716

# ALC row 76
$ws.Range("H76").Value = 6176073
$ws.Range("I76").Value = 6176073
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 6176073
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -6175758

# ALC row 79
$ws.Range("H79").Value = 6176073
$ws.Range("I79").Value = 6176073
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 6176073
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -6174981

# ALC row 100
$ws.Range("H100").Value = 13890575
$ws.Range("I100").Value = 17544776
$ws.Range("J100").Value = 4605.2
$ws.Range("K100").Value = 17544776
$ws.Range("L100").Value = 4605.2
$ws.Range("M100").Value = -17544235
$ws.Range("N100").Value = -5687.2

# ALC row 116
$ws.Range("H116").Value = 4805.1
$ws.Range("I116").Value = 6683.8335
$ws.Range("J116").Value = 1987
$ws.Range("K116").Value = 6683.8335
$ws.Range("L116").Value = 1987
$ws.Range("M116").Value = -3241.8335
$ws.Range("N116").Value = -8871

# ALC row 135
$ws.Range("H135").Value = 1506.7931
$ws.Range("I135").Value = 1721.7916
$ws.Range("J135").Value = 474.8
$ws.Range("K135").Value = 15496.1244
$ws.Range("L135").Value = 4273.2
$ws.Range("M135").Value = -12961.1244

# ALC row 137
$ws.Range("H137").Value = 1066.5714
$ws.Range("I137").Value = 1001.64703
$ws.Range("J137").Value = 1127.8889
$ws.Range("K137").Value = 3004.94109
$ws.Range("L137").Value = 3383.6667
$ws.Range("M137").Value = -454.9410899999998

# ALC row 139
$ws.Range("H139").Value = 77600
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 77600
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 77600
$ws.Range("N139").Value = -87880

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 453369.3
$ws.Range("I32").Value = 3981.4546
$ws.Range("J32").Value = 3199628.5
$ws.Range("K32").Value = 3981.4546
$ws.Range("L32").Value = 3199628.5
$ws.Range("M32").Value = -3694.4546
$ws.Range("N32").Value = -3200202.5

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2622.9285
$ws.Range("I107").Value = 2610.7
$ws.Range("J107").Value = 2653.5
$ws.Range("K107").Value = 2610.7
$ws.Range("L107").Value = 2653.5
$ws.Range("M107").Value = -690.6999999999998
$ws.Range("N107").Value = -6493.5

# BSM row 134
$ws.Range("H134").Value = 2091.1555
$ws.Range("I134").Value = 2000.8649
$ws.Range("J134").Value = 2508.75
$ws.Range("K134").Value = 6002.5947
$ws.Range("L134").Value = 7526.25
$ws.Range("M134").Value = -3467.5947
$ws.Range("N134").Value = -12596.25

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2498

# CRP row 65
$ws.Range("H65").Value = 2498

# CRP row 99
$ws.Range("H99").Value = 426601.06
$ws.Range("I99").Value = 473822.28
$ws.Range("J99").Value = 1610
$ws.Range("K99").Value = 473822.28
$ws.Range("L99").Value = 1610
$ws.Range("M99").Value = -472324.28
$ws.Range("N99").Value = -4606

# CRP row 126
$ws.Range("H126").Value = 426601.06
$ws.Range("I126").Value = 473822.28
$ws.Range("J126").Value = 1610
$ws.Range("K126").Value = 1421466.84
$ws.Range("L126").Value = 4830
$ws.Range("M126").Value = -1418996.84
$ws.Range("N126").Value = -9770

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 883.4167
$ws.Range("I68").Value = 650
$ws.Range("J68").Value = 930.1
$ws.Range("K68").Value = 1950
$ws.Range("L68").Value = 2790.3
$ws.Range("M68").Value = -1139
$ws.Range("N68").Value = -4412.3

# CUL row 71
$ws.Range("H71").Value = 883.4167
$ws.Range("I71").Value = 650
$ws.Range("J71").Value = 930.1
$ws.Range("K71").Value = 5850
$ws.Range("L71").Value = 8370.9
$ws.Range("M71").Value = -1794
$ws.Range("N71").Value = -16482.9

# GSM row 64
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# GSM row 67
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# GSM row 70
$ws.Range("H70").Value = 17052674
$ws.Range("I70").Value = 35163196
$ws.Range("J70").Value = 7479.4116
$ws.Range("K70").Value = 35163196
$ws.Range("L70").Value = 7479.4116
$ws.Range("M70").Value = -35162926
$ws.Range("N70").Value = -8019.4116

# GSM row 73
$ws.Range("H73").Value = 17052674
$ws.Range("I73").Value = 35163196
$ws.Range("J73").Value = 7479.4116
$ws.Range("K73").Value = 35163196
$ws.Range("L73").Value = 7479.4116
$ws.Range("M73").Value = -35162260
$ws.Range("N73").Value = -9351.411599999999

# GSM row 80
$ws.Range("H80").Value = 4323.4443
$ws.Range("I80").Value = 3501.6667
$ws.Range("J80").Value = 4734.3335
$ws.Range("K80").Value = 3501.6667
$ws.Range("L80").Value = 4734.3335
$ws.Range("M80").Value = -2503.6667
$ws.Range("N80").Value = -6730.3335

# GSM row 83
$ws.Range("H83").Value = 4323.4443
$ws.Range("I83").Value = 3501.6667
$ws.Range("J83").Value = 4734.3335
$ws.Range("K83").Value = 17508.3335
$ws.Range("L83").Value = 23671.6675
$ws.Range("M83").Value = -12516.3335
$ws.Range("N83").Value = -33655.6675

# GSM row 107
$ws.Range("H107").Value = 533
$ws.Range("I107").Value = 533
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 533
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1387
$ws.Range("N107").ClearContents()

# GSM row 122
$ws.Range("H122").Value = 2131.889
$ws.Range("I122").Value = 2026.7142
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 6080.142599999999
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -3630.142599999999

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2408.8333
$ws.Range("I68").Value = 2100
$ws.Range("J68").Value = 3335.3333
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 3335.3333
$ws.Range("M68").Value = -1351
$ws.Range("N68").Value = -4833.3333

# LTW row 71
$ws.Range("H71").Value = 2408.8333
$ws.Range("I71").Value = 2100
$ws.Range("J71").Value = 3335.3333
$ws.Range("K71").Value = 10500
$ws.Range("L71").Value = 16676.6665
$ws.Range("M71").Value = -6756
$ws.Range("N71").Value = -24164.6665

# LTW row 74
$ws.Range("H74").Value = 15800
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 15800
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 15800
$ws.Range("N74").Value = -17796

# LTW row 77
$ws.Range("H77").Value = 15800
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 15800
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 47400
$ws.Range("N77").Value = -57384

# WVR row 105
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 19995
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 19995
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 19995
$ws.Range("N105").Value = -26983
